# Apply cryptos list price/volume update (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving exact characters (no numeric/date autoconversion)
# and without leaving the cell tagged with a non-default "Text" number format/style.
function Set-CellText {
    param($Cell, [string]$Text)
    $Cell.Value = "'" + $Text
    $Cell.Style = "Normal"
}

Set-CellText $ws.Range("D2") "69.153.58"
Set-CellText $ws.Range("E2") "  +0.62%  "

Set-CellText $ws.Range("D3") "3.770.38"
Set-CellText $ws.Range("E3") "  +0.39%  "

Set-CellText $ws.Range("E4") "  +0.01%  "

Set-CellText $ws.Range("D5") "602.03"
Set-CellText $ws.Range("E5") "  +0.07%  "

Set-CellText $ws.Range("D6") "166.09"
Set-CellText $ws.Range("E6") "  -1.78%  "

Set-CellText $ws.Range("D7") "3.769.00"
Set-CellText $ws.Range("E7") "  +0.38%  "

Set-CellText $ws.Range("E8") "  +0.03%  "

Set-CellText $ws.Range("E9") "  +0.42%  "

Set-CellText $ws.Range("E10") "  +4.35%  "

Set-CellText $ws.Range("D11") "6.35"
Set-CellText $ws.Range("E11") "  +0.18%  "

Set-CellText $ws.Range("E12") "  -0.48%  "

Set-CellText $ws.Range("D13") "37.67"
Set-CellText $ws.Range("E13") "  -1.67%  "

Set-CellText $ws.Range("E14") "  +0.18%  "

Set-CellText $ws.Range("D15") "4.401.43"
Set-CellText $ws.Range("E15") "  +0.42%  "

Set-CellText $ws.Range("D16") "3.784.83"
Set-CellText $ws.Range("E16") "  +0.73%  "

Set-CellText $ws.Range("D17") "69.279.16"
Set-CellText $ws.Range("E17") "  +0.79%  "

Set-CellText $ws.Range("E18") "  +1.64%  "

Set-CellText $ws.Range("D19") "17.64"
Set-CellText $ws.Range("E19") "  +3.17%  "

Set-CellText $ws.Range("E20") "  -1.06%  "

Set-CellText $ws.Range("D21") "11.29"
Set-CellText $ws.Range("E21") "  +3.67%  "

Set-CellText $ws.Range("D22") "491.97"
Set-CellText $ws.Range("E22") "  -0.70%  "

Set-CellText $ws.Range("E23") "  -0.58%  "

Set-CellText $ws.Range("E24") "  -1.92%  "

Set-CellText $ws.Range("D25") "84.82"
Set-CellText $ws.Range("E25") "  -0.52%  "

Set-CellText $ws.Range("E26") "  -2.69%  "

Set-CellText $ws.Range("E27") "  -0.78%  "

Set-CellText $ws.Range("D28") "10.09"
Set-CellText $ws.Range("E28") "  -1.77%  "

Set-CellText $ws.Range("E29") "  +0.07%  "

Set-CellText $ws.Range("E30") "  -0.41%  "

Set-CellText $ws.Range("D31") "8.13"
Set-CellText $ws.Range("E31") "  +2.76%  "

Set-CellText $ws.Range("E32") "  -4.06%  "

Set-CellText $ws.Range("D33") "31.84"
Set-CellText $ws.Range("E33") "  -0.24%  "

Set-CellText $ws.Range("D34") "3.914.18"
Set-CellText $ws.Range("E34") "  +0.34%  "

Set-CellText $ws.Range("D35") "3.723.87"
Set-CellText $ws.Range("E35") "  +0.93%  "

Set-CellText $ws.Range("E36") "  -0.79%  "

Set-CellText $ws.Range("D37") "5.95"
Set-CellText $ws.Range("E37") "  +1.46%  "

Set-CellText $ws.Range("E38") "  -0.20%  "

Set-CellText $ws.Range("E39") "  +4.00%  "

Set-CellText $ws.Range("E40") "  +0.06%  "

Set-CellText $ws.Range("E41") "  +5.44%  "

Set-CellText $ws.Range("E42") "  +0.24%  "

Set-CellText $ws.Range("B43") "OKB"
Set-CellText $ws.Range("C43") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-CellText $ws.Range("D43") "48.54"
Set-CellText $ws.Range("E43") "  -0.75%  "

Set-CellText $ws.Range("B44") "Bittensor"
Set-CellText $ws.Range("C44") "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-CellText $ws.Range("D44") "426.26"
Set-CellText $ws.Range("E44") "  -3.34%  "

Set-CellText $ws.Range("B45") "Stacks"
Set-CellText $ws.Range("C45") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-CellText $ws.Range("D45") "1.99"
Set-CellText $ws.Range("E45") "  +0.60%  "

Set-CellText $ws.Range("D46") "8.43"
Set-CellText $ws.Range("E46") "  -0.59%  "

Set-CellText $ws.Range("E47") "  +0.05%  "

Set-CellText $ws.Range("B48") "Monero"
Set-CellText $ws.Range("C48") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-CellText $ws.Range("D48") "142.16"
Set-CellText $ws.Range("E48") "  +0.42%  "

Set-CellText $ws.Range("B49") "Arweave"
Set-CellText $ws.Range("C49") "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-CellText $ws.Range("D49") "40.02"
Set-CellText $ws.Range("E49") "  -1.09%  "

Set-CellText $ws.Range("D50") "2.806.91"
Set-CellText $ws.Range("E50") "  +0.47%  "

Set-CellText $ws.Range("E51") "  +8.02%  "
